$wb = $excel.ActiveWorkbook

# --- Update status text: "Ready for handoff" -> "In Translation" ---
$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsDeDe = $wb.Worksheets.Item("de-de")

$wsOverview.Range("E2").Value = "In Translation"
$wsOverview.Range("F2").Value = "In Translation"
$wsZhCn.Range("C2").Value = "In Translation"
$wsDeDe.Range("C2").Value = "In Translation"

# --- Narrow the "Status" related columns ---
# (COM ColumnWidth is quantized to 1/6-character increments by the engine;
#  12.5 is the input that lands on the closest representable stored width
#  to the target 13.4101845877511, i.e. 13.333333333333334.)
$wsOverview.Columns(5).ColumnWidth = 12.5
$wsOverview.Columns(6).ColumnWidth = 12.5
$wsZhCn.Columns(3).ColumnWidth = 12.5
$wsDeDe.Columns(3).ColumnWidth = 12.5
